$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "302.91"
    "E2" = "3.77%"
    "D3" = "34.82"
    "E3" = "12.25%"
    "D4" = "5.157"
    "E4" = "4.78%"
    "D5" = "0.07753"
    "E5" = "4.80%"
    "D6" = "2.379"
    "E6" = "7.80%"
    "D7" = "8.012"
    "E7" = "3.61%"
    "D8" = "3.947"
    "E8" = "5.24%"
    "D9" = "0.9295"
    "E9" = "1.89%"
    "D10" = "0.09859"
    "E10" = "10.51%"
    "D11" = "0.1797"
    "E11" = "6.48%"
    "D12" = "0.08640"
    "E12" = "4.54%"
    "D13" = "0.03318"
    "E13" = "6.74%"
    "D14" = "0.09892"
    "E14" = "-0.99%"
    "D15" = "0.001501"
    "E15" = "0.16%"
    "D16" = "0.005775"
    "E16" = "-1.32%"
    "E17" = "-1.03%"
    "D18" = "2.134"
    "E18" = "2.23%"
    "D19" = "0.3368"
    "E19" = "1.23%"
    "E20" = "2.73%"
    "D21" = "4.323"
    "E21" = "8.53%"
    "E22" = "5.08%"
    "D23" = "0.04576"
    "E23" = "0.40%"
    "D24" = "0.001217"
    "E24" = "0.20%"
    "E25" = "-2.58%"
    "D26" = "0.0001300"
    "E26" = "-0.20%"
    "E27" = "-0.27%"
    "D39" = "0.01789"
    "E39" = "12.38%"
    "D40" = "0.04795"
    "E40" = "7.37%"
    "D41" = "0.007744"
    "E41" = "5.11%"
    "D42" = "0.1411"
    "E42" = "6.34%"
    "D43" = "0.007101"
    "E43" = "-26.63%"
    "D44" = "0.002140"
    "E44" = "-8.34%"
    "D45" = "0.009178"
    "E45" = "0.37%"
    "D46" = "0.00006121"
    "E46" = "0.25%"
    "E47" = "-0.16%"
    "E48" = "44.94%"
    "D49" = "0.002000"
    "E49" = "-0.21%"
    "D50" = "0.00002100"
    "E50" = "-0.16%"
    "D51" = "0.0002000"
    "E51" = "-0.16%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
    $range.Style = "Normal"
}
